$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSheetFor Saving MatrixData")
$ws.Range("B1").Replace("Invalid cast from 'System.Double' to 'System.Nullable`1[[System.Double, mscorlib, Version=4.0.0.0, Culture=neutral, PublicKeyToken=b77a5c561934e089]]'.", "A:8")
Write-Host $ws.Range("B1").Text
Write-Host $ws.Range("B1").Formula
